$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.341.12"

$ws.Range("D3").Value = "1.879.71"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "0.7103"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").Value = "242.47"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.08027"
$ws.Range("E8").Value = "  +3.19%  "

$ws.Range("D9").Value = "0.3141"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").Value = "25.12"

$ws.Range("D11").Value = "0.08326"

$ws.Range("D12").Value = "1.897.31"
$ws.Range("E12").Value = "  +1.37%  "

$ws.Range("D13").Value = "5.269"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").Value = "94.60"
$ws.Range("E14").Value = "  +3.77%  "

$ws.Range("D15").Value = "0.7181"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").Value = "6.352"
$ws.Range("E16").Value = "  +5.19%  "

$ws.Range("D17").Value = "0.000008725"
$ws.Range("E17").Value = "  +5.75%  "

$ws.Range("D18").Value = "29.363.29"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("D19").Value = "243.05"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Value = "2.150.64"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").Value = "13.32"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D23").Value = "7.866"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "0.1573"
$ws.Range("E25").Value = "  -2.32%  "

$ws.Range("D26").Value = "163.64"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "9.076"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  -0.49%  "

$ws.Range("D30").Value = "4.438"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("E31").Value = "  +1.20%  "

$ws.Range("E32").Value = "  -6.21%  "

$ws.Range("D33").Value = "0.05394"

$ws.Range("D34").Value = "1.942"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").Value = "0.7769"
$ws.Range("E35").Value = "  +3.97%  "

$ws.Range("D36").Value = "1.179"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").Value = "0.01886"
$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").Value = "1.270.15"
$ws.Range("E39").Value = "  +5.37%  "

$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").Value = "6.554"
$ws.Range("E41").Value = "  +1.61%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.9213"
$ws.Range("E42").Value = "  +3.74%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "113.80"
$ws.Range("E43").Value = "  +4.43%  "

$ws.Range("E44").Value = "  +2.49%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").Value = "2.050.13"

$ws.Range("E47").Value = "  +3.75%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.805"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.5220"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").Value = "9.555"
$ws.Range("E50").Value = "  +1.80%  "

$ws.Range("D51").Value = "0.4379"
$ws.Range("E51").Value = "  +1.43%  "

Write-Output "done"